{"js": "// Update stack trace line numbers / frames inside the big error dump\n// paragraph (M2Doc test fixture) to reflect the library bump described\n// in the commit message (\"Moving from 2.0.2 to 2.0.3\").\nconst body = context.document.body;\n\n// 1) A handful of single line-number tweaks; every search string below\n//    is unique within the document, so a straight search+replace is safe.\nconst simpleReplacements = [\n  [\n    \"JavaMethodService.internalInvoke(JavaMethodService.java:163)\",\n    \"JavaMethodService.internalInvoke(JavaMethodService.java:162)\",\n  ],\n  [\n    \"AbstractService.invoke(AbstractService.java:136)\",\n    \"AbstractService.invoke(AbstractService.java:135)\",\n  ],\n  [\n    \"EvaluationServices.call(EvaluationServices.java:168)\",\n    \"EvaluationServices.call(EvaluationServices.java:172)\",\n  ],\n  [\n    \"EvaluationServices.callOrApply(EvaluationServices.java:204)\",\n    \"EvaluationServices.callOrApply(EvaluationServices.java:208)\",\n  ],\n  [\n    \"AstEvaluator.caseCall(AstEvaluator.java:192)\",\n    \"AstEvaluator.caseCall(AstEvaluator.java:189)\",\n  ],\n  [\n    \"AstSwitch.doSwitch(AstSwitch.java:118)\",\n    \"AstSwitch.doSwitch(AstSwitch.java:119)\",\n  ],\n  [\n    \"AstEvaluator.eval(AstEvaluator.java:112)\",\n    \"AstEvaluator.eval(AstEvaluator.java:109)\",\n  ],\n  [\n    \"GeneratedMethodAccessor74.invoke(Unknown Source)\",\n    \"GeneratedMethodAccessor73.invoke(Unknown Source)\",\n  ],\n];\n\nfor (const [oldText, newText] of simpleReplacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 2) The tail of the stack trace (Maven/Tycho/Equinox launcher frames)\n//    is replaced wholesale with the JDT/Eclipse JUnit runner frames.\nconst startAnchorText =\n  \"org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\";\nconst endAnchorText =\n  \"org.eclipse.equinox.launcher.Main.main(Main.java:1471)\";\n\nconst startResults = body.search(startAnchorText, { matchCase: true });\nconst endResults = body.search(endAnchorText, { matchCase: true });\ncontext.load(startResults, \"text\");\ncontext.load(endResults, \"text\");\nawait context.sync();\n\nconst blockRange = startResults.items[0].expandTo(endResults.items[0]);\n\nconst newTail = [\n  \"org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\",\n].join(\"\\n\");\n\nblockRange.insertText(newTail, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update stack trace line numbers / frames inside the big error dump\n# paragraph (M2Doc test fixture) to reflect the library bump described\n# in the commit message (\"Moving from 2.0.2 to 2.0.3\").\n$d = $word.ActiveDocument\n\n# 1) A handful of single line-number tweaks; every search string below\n#    is unique within the document, so a straight Find/Replace is safe.\n$simplePairs = @(\n    @(\"JavaMethodService.internalInvoke(JavaMethodService.java:163)\", \"JavaMethodService.internalInvoke(JavaMethodService.java:162)\"),\n    @(\"AbstractService.invoke(AbstractService.java:136)\", \"AbstractService.invoke(AbstractService.java:135)\"),\n    @(\"EvaluationServices.call(EvaluationServices.java:168)\", \"EvaluationServices.call(EvaluationServices.java:172)\"),\n    @(\"EvaluationServices.callOrApply(EvaluationServices.java:204)\", \"EvaluationServices.callOrApply(EvaluationServices.java:208)\"),\n    @(\"AstEvaluator.caseCall(AstEvaluator.java:192)\", \"AstEvaluator.caseCall(AstEvaluator.java:189)\"),\n    @(\"AstSwitch.doSwitch(AstSwitch.java:118)\", \"AstSwitch.doSwitch(AstSwitch.java:119)\"),\n    @(\"AstEvaluator.eval(AstEvaluator.java:112)\", \"AstEvaluator.eval(AstEvaluator.java:109)\"),\n    @(\"GeneratedMethodAccessor74.invoke(Unknown Source)\", \"GeneratedMethodAccessor73.invoke(Unknown Source)\")\n)\n\nforeach ($pair in $simplePairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n# 2) The tail of the stack trace (Maven/Tycho/Equinox launcher frames)\n#    is replaced wholesale with the JDT/Eclipse JUnit runner frames.\n$r1 = $d.Range()\n$f1 = $r1.Find\n$f1.Text = \"org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\"\n$f1.Execute() | Out-Null\n\n$r2 = $d.Range()\n$f2 = $r2.Find\n$f2.Text = \"org.eclipse.equinox.launcher.Main.main(Main.java:1471)\"\n$f2.Execute() | Out-Null\n\n$blockRange = $d.Range($r1.Start, $r2.End)\n\n$newTailLines = @(\n    \"org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\",\n    \"at org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\",\n    \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)\",\n    \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)\",\n    \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)\",\n    \"at org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)\"\n)\n$newTail = $newTailLines[0] + \"`n`t\" + ($newTailLines[1..($newTailLines.Length - 1)] -join \"`n`t\")\n\n$blockRange.Text = $newTail\n"}
